$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in W1, matching style/formatting of existing headers (e.g. V1)
$ws.Range("W1").Value = "param_E_pv3_solar"
$ws.Range("V1").Copy()
$ws.Range("W1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New values for columns U, V (overwritten) and W (newly added) for rows 2-16
$data = @{
    2  = @(0.12, 42.35, 52.40552777777779)
    3  = @(0.12, 69.46500000000002, 65.41700000000003)
    4  = @(0.12, 66.00000000000004, 34.73616666666666)
    5  = @(0.12, 0, 0)
    6  = @(0.12, 0, 111.9891666666667)
    7  = @(0.12, 49.93404166666668, 33)
    8  = @(0.12, 33.37400000000001, 0)
    9  = @(0.12, 157.15425, 0)
    10 = @(0.12, 0, 0)
    11 = @(0.12, 0, 0)
    12 = @(0.12, 0, 0)
    13 = @(0.12, 0, 0)
    14 = @(0.12, 0, 0)
    15 = @(0.12, 0, 0)
    16 = @(0.12, 0, 0)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("U$row").Value = $vals[0]
    $ws.Range("V$row").Value = $vals[1]
    $ws.Range("W$row").Value = $vals[2]
}
